$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("ID in folder"): multiply every ID by 10 to make space for new IDs between them
$ws.Range("C2").Value = 10
$ws.Range("C3").Value = 20
$ws.Range("C4").Value = 30
$ws.Range("C5").Value = 40
$ws.Range("C6").Value = 50
$ws.Range("C7").Value = 60
$ws.Range("C8").Value = 70
$ws.Range("C9").Value = 80
$ws.Range("C10").Value = 90
$ws.Range("C11").Value = 100
$ws.Range("C12").Value = 110
$ws.Range("C13").Value = 120
$ws.Range("C14").Value = 130
$ws.Range("C15").Value = 140
$ws.Range("C16").Value = 150
$ws.Range("C17").Value = 160
$ws.Range("C18").Value = 170
$ws.Range("C19").Value = 180
$ws.Range("C20").Value = 190
$ws.Range("C21").Value = 200
$ws.Range("C22").Value = 210
$ws.Range("C23").Value = 220
$ws.Range("C24").Value = 230
$ws.Range("C25").Value = 240
$ws.Range("C26").Value = 250
$ws.Range("C27").Value = 260
$ws.Range("C28").Value = 270
$ws.Range("C29").Value = 280
$ws.Range("C30").Value = 290
$ws.Range("C31").Value = 300
$ws.Range("C32").Value = 310
$ws.Range("C33").Value = 320
$ws.Range("C34").Value = 330
$ws.Range("C35").Value = 340
$ws.Range("C36").Value = 350
$ws.Range("C37").Value = 360
$ws.Range("C38").Value = 370
$ws.Range("C39").Value = 380
$ws.Range("C40").Value = 390
$ws.Range("C41").Value = 400
$ws.Range("C42").Value = 410

# Column D ("ID in Thesis"): clear out old thesis-ID numbering (no longer aligned with new IDs)
$ws.Range("D3").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("D18").Value = ""
$ws.Range("D20").Value = ""
$ws.Range("D23").Value = ""
$ws.Range("D24").Value = ""
$ws.Range("D26").Value = ""
$ws.Range("D29").Value = ""
$ws.Range("D31").Value = ""
$ws.Range("D37").Value = ""

# Column L ("Notes"): update references to other TC numbers to match the new *10 numbering
$ws.Range("L9").Value = 'Results are together with 70'
$ws.Range("L11").Value = 'Removed as too similar to 50'
$ws.Range("L18").Value = 'Combined with TC180 and 210'
$ws.Range("L19").Value = 'Done in TC170'
$ws.Range("L21").Value = 'Done in TC190'
$ws.Range("L22").Value = 'Done in TC170'
$ws.Range("L25").Value = 'Covered in TC 230'
$ws.Range("L27").Value = 'Covered in TC250'
$ws.Range("L32").Value = 'Could take a long time and not be worth it. Similar to TC 300'
$ws.Range("L33").Value = 'Covered in TC 300'
$ws.Range("L34").Value = 'Covered in TC 310'
$ws.Range("L35").Value = 'Covered in TC 300'
$ws.Range("L36").Value = 'Covered in TC 310'

# Restore the selection/scroll state
$ws.Range("F7").Select()
